$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was added on top of the existing "Platano" price
# history. This pushes every existing record (rows 99-183) down by one row
# (the oldest record, which was in row 183, ends up in the newly created
# row 184), and the brand-new record is written into row 99.

# Insert a new row at row 99; Excel shifts rows 99-183 down to 100-184.
$ws.Rows.Item(99).Insert()

# Populate the new row 99 with the new daily record.
$ws.Cells.Item(99, 1).Value = 1
$ws.Cells.Item(99, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(99, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(99, 4).Value = 44586
$ws.Cells.Item(99, 5).Value = 15
$ws.Cells.Item(99, 6).Value = "Fruta"
$ws.Cells.Item(99, 7).Value = 100108
$ws.Cells.Item(99, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(99, 9).Value = 100108006
$ws.Cells.Item(99, 10).Value = "Plátano"
$ws.Cells.Item(99, 11).Value = "Sin especificar"
$ws.Cells.Item(99, 12).Value = "Pintón"
$ws.Cells.Item(99, 13).Value = 120
$ws.Cells.Item(99, 14).Value = 19000
$ws.Cells.Item(99, 15).Value = 20000
$ws.Cells.Item(99, 16).Value = 19500
$ws.Cells.Item(99, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(99, 18).Value = "Ecuador"
$ws.Cells.Item(99, 19).Value = 975
$ws.Cells.Item(99, 20).Value = 20
